$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the stray empty styled cell left over in row 18
$ws.Rows("18").Delete() | Out-Null

# Add the new logged entry as row 13
$ws.Range("B13").Value = "Programmazione"
$ws.Range("C13").Value = "Unity e Ink"
$ws.Range("D13").Value = 0.041666666666666664
$ws.Range("E13").Value = "Funzioni di cambio luogo + background base"

# Move the active selection to where the next entry would go
$ws.Range("E21").Select() | Out-Null
